$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-5 with new values
$data = @(
    @(1, 5, 3, 1, 5, -4, 2, 54, 5),
    @(2, 6, 2, 1, 3, -5, 1, 65, 5),
    @(3, 6, 4, 5, 9, -1, 5, 21, 5),
    @(4, 5, 1, 2, 4, -3, 3, 43, 5),
    @(5, 8, 4, 6, 8, -2, 4, 32, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
}

$ws.Range("I1").Select()
